# Refresh computed market-price / leve-profit columns (H:N) in the Tiamat Profits workbook.
# Values mirror a scheduled market-board data pull; only numeric outputs change, no
# structural/layout edits. Two rows (ALC!86 and ALC!89) drop their now-redundant HQ-profit
# (column N) cell once NQ/HQ pricing collapses to a single figure.
$wb = $excel.ActiveWorkbook

# ============ Sheet: ALC ============
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value2 = 203.33333
$ws.Range("I11").Value2 = 203.33333
$ws.Range("K11").Value2 = 203.33333
$ws.Range("M11").Value2 = -63.33332999999999
# Row 12
$ws.Range("H12").Value2 = 180.73077
$ws.Range("I12").Value2 = 187.47058
$ws.Range("J12").Value2 = 168
$ws.Range("K12").Value2 = 187.47058
$ws.Range("L12").Value2 = 168
$ws.Range("M12").Value2 = -17.47058000000001
$ws.Range("N12").Value2 = -508
# Row 18
$ws.Range("H18").Value2 = 283.92856
$ws.Range("I18").Value2 = 147.91667
$ws.Range("J18").Value2 = 1100
$ws.Range("K18").Value2 = 147.91667
$ws.Range("L18").Value2 = 1100
$ws.Range("M18").Value2 = 136.08333
$ws.Range("N18").Value2 = -1668
# Row 33
$ws.Range("H33").Value2 = 600.26086
$ws.Range("I33").Value2 = 69.125
$ws.Range("J33").Value2 = 1814.2858
$ws.Range("K33").Value2 = 69.125
$ws.Range("L33").Value2 = 1814.2858
$ws.Range("M33").Value2 = 159.875
$ws.Range("N33").Value2 = -2272.2858
# Row 41
$ws.Range("H41").Value2 = 362.1
$ws.Range("I41").Value2 = 320.33334
$ws.Range("J41").Value2 = 424.75
$ws.Range("K41").Value2 = 320.33334
$ws.Range("L41").Value2 = 424.75
$ws.Range("M41").Value2 = 119.66666
$ws.Range("N41").Value2 = -1304.75
# Row 53
$ws.Range("H53").Value2 = 421.45456
$ws.Range("I53").Value2 = 466.66666
$ws.Range("J53").Value2 = 404.5
$ws.Range("K53").Value2 = 466.66666
$ws.Range("L53").Value2 = 404.5
$ws.Range("M53").Value2 = 170.33334
$ws.Range("N53").Value2 = -1678.5
# Row 70
$ws.Range("H70").Value2 = 2826.4
$ws.Range("I70").Value2 = 1500
$ws.Range("J70").Value2 = 2973.7778
$ws.Range("K70").Value2 = 4500
$ws.Range("L70").Value2 = 8921.3334
$ws.Range("M70").Value2 = -4230
$ws.Range("N70").Value2 = -9461.3334
# Row 73
$ws.Range("H73").Value2 = 2826.4
$ws.Range("I73").Value2 = 1500
$ws.Range("J73").Value2 = 2973.7778
$ws.Range("K73").Value2 = 4500
$ws.Range("L73").Value2 = 8921.3334
$ws.Range("M73").Value2 = -3564
$ws.Range("N73").Value2 = -10793.3334
# Row 76
$ws.Range("H76").Value2 = 23259128
$ws.Range("I76").Value2 = 27030126
$ws.Range("J76").Value2 = 4635.5
$ws.Range("K76").Value2 = 27030126
$ws.Range("L76").Value2 = 4635.5
$ws.Range("M76").Value2 = -27029811
$ws.Range("N76").Value2 = -5265.5
# Row 79
$ws.Range("H79").Value2 = 23259128
$ws.Range("I79").Value2 = 27030126
$ws.Range("J79").Value2 = 4635.5
$ws.Range("K79").Value2 = 27030126
$ws.Range("L79").Value2 = 4635.5
$ws.Range("M79").Value2 = -27029034
$ws.Range("N79").Value2 = -6819.5
# Row 86
$ws.Range("H86").Value2 = 1219.8667
$ws.Range("I86").Value2 = 1219.8667
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 1219.8667
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = -96.86670000000004
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value2 = 1219.8667
$ws.Range("I89").Value2 = 1219.8667
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 6099.333500000001
$ws.Range("L89").Value2 = 0
$ws.Range("M89").Value2 = -483.3335000000006
$ws.Range("N89").ClearContents()
# Row 92
$ws.Range("H92").Value2 = 29412462
$ws.Range("I92").Value2 = 62500670
$ws.Range("J92").Value2 = 720.55554
$ws.Range("K92").Value2 = 62500670
$ws.Range("L92").Value2 = 720.55554
$ws.Range("M92").Value2 = -62499422
$ws.Range("N92").Value2 = -3216.55554
# Row 94
$ws.Range("H94").Value2 = 8890.4
$ws.Range("I94").Value2 = 8333
$ws.Range("J94").Value2 = 9726.5
$ws.Range("K94").Value2 = 8333
$ws.Range("L94").Value2 = 9726.5
$ws.Range("M94").Value2 = -7882
$ws.Range("N94").Value2 = -10628.5
# Row 98
$ws.Range("H98").Value2 = 1067.2
$ws.Range("I98").Value2 = 1086
$ws.Range("J98").Value2 = 898
$ws.Range("K98").Value2 = 1086
$ws.Range("L98").Value2 = 898
$ws.Range("M98").Value2 = 412
$ws.Range("N98").Value2 = -3894
# Row 107
$ws.Range("H107").Value2 = 571.1070999999999
$ws.Range("I107").Value2 = 582.9583
$ws.Range("K107").Value2 = 582.9583
$ws.Range("M107").Value2 = 1337.0417
# Row 113
$ws.Range("H113").Value2 = 1873.3462
$ws.Range("I113").Value2 = 1800.2778
$ws.Range("J113").Value2 = 2037.75
$ws.Range("K113").Value2 = 1800.2778
$ws.Range("L113").Value2 = 2037.75
$ws.Range("M113").Value2 = 1453.7222
$ws.Range("N113").Value2 = -8545.75
# Row 116
$ws.Range("H116").Value2 = 3608.6667
$ws.Range("I116").Value2 = 2923
$ws.Range("K116").Value2 = 2923
$ws.Range("M116").Value2 = 519
# Row 122
$ws.Range("H122").Value2 = 1067.2
$ws.Range("I122").Value2 = 1086
$ws.Range("J122").Value2 = 898
$ws.Range("K122").Value2 = 3258
$ws.Range("L122").Value2 = 2694
$ws.Range("M122").Value2 = -808
$ws.Range("N122").Value2 = -7594
# Row 129
$ws.Range("H129").Value2 = 946.6070999999999
$ws.Range("J129").Value2 = 1000.73914
$ws.Range("L129").Value2 = 3002.21742
$ws.Range("N129").Value2 = -13002.21742
# Row 135
$ws.Range("H135").Value2 = 11364309
$ws.Range("I135").Value2 = 456.79413
$ws.Range("J135").Value2 = 50001410
$ws.Range("K135").Value2 = 4111.14717
$ws.Range("L135").Value2 = 450012690
$ws.Range("M135").Value2 = -1576.14717
$ws.Range("N135").Value2 = -450017760
# Row 137
$ws.Range("H137").Value2 = 6223.0527
$ws.Range("I137").Value2 = 697.2
$ws.Range("K137").Value2 = 2091.6
$ws.Range("M137").Value2 = 458.3999999999996

# ============ Sheet: ARM ============
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value2 = 4303.643
$ws.Range("I2").Value2 = 733
$ws.Range("K2").Value2 = 733
$ws.Range("M2").Value2 = -620
# Row 45
$ws.Range("H45").Value2 = 851.4167
$ws.Range("I45").Value2 = 815.38464
$ws.Range("J45").Value2 = 894
$ws.Range("K45").Value2 = 815.38464
$ws.Range("L45").Value2 = 894
$ws.Range("M45").Value2 = -438.38464
$ws.Range("N45").Value2 = -1648
# Row 61
$ws.Range("H61").Value2 = 2681.9707
$ws.Range("I61").Value2 = 2166.1853
$ws.Range("J61").Value2 = 4671.4287
$ws.Range("K61").Value2 = 2166.1853
$ws.Range("L61").Value2 = 4671.4287
$ws.Range("M61").Value2 = -1954.1853
$ws.Range("N61").Value2 = -5095.4287
# Row 97
$ws.Range("H97").Value2 = 1102.7273
$ws.Range("I97").Value2 = 651.0526
$ws.Range("J97").Value2 = 3963.3333
$ws.Range("K97").Value2 = 651.0526
$ws.Range("L97").Value2 = 3963.3333
$ws.Range("M97").Value2 = -155.0526
$ws.Range("N97").Value2 = -4955.3333
# Row 110
$ws.Range("H110").Value2 = 1732.8572
$ws.Range("I110").Value2 = 1667.7646
$ws.Range("J110").Value2 = 2009.5
$ws.Range("K110").Value2 = 1667.7646
$ws.Range("L110").Value2 = 2009.5
$ws.Range("M110").Value2 = 377.2354
$ws.Range("N110").Value2 = -6099.5
# Row 116
$ws.Range("H116").Value2 = 4303.643
$ws.Range("I116").Value2 = 733
$ws.Range("K116").Value2 = 733
$ws.Range("M116").Value2 = 1561
# Row 122
$ws.Range("H122").Value2 = 1485.9
$ws.Range("I122").Value2 = 1228
$ws.Range("J122").Value2 = 1657.8334
$ws.Range("K122").Value2 = 3684
$ws.Range("L122").Value2 = 4973.5002
$ws.Range("M122").Value2 = -1234
$ws.Range("N122").Value2 = -9873.5002
# Row 136
$ws.Range("H136").Value2 = 2681.9707
$ws.Range("I136").Value2 = 2166.1853
$ws.Range("J136").Value2 = 4671.4287
$ws.Range("K136").Value2 = 6498.5559
$ws.Range("L136").Value2 = 14014.2861
$ws.Range("M136").Value2 = -3948.5559
$ws.Range("N136").Value2 = -19114.2861

# ============ Sheet: BSM ============
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value2 = 4303.643
$ws.Range("I3").Value2 = 733
$ws.Range("K3").Value2 = 733
$ws.Range("M3").Value2 = -619
# Row 20
$ws.Range("H20").Value2 = 1359.1389
$ws.Range("I20").Value2 = 1410
$ws.Range("K20").Value2 = 1410
$ws.Range("M20").Value2 = -1163
# Row 134
$ws.Range("H134").Value2 = 22246278
$ws.Range("I134").Value2 = 1657.2858
$ws.Range("J134").Value2 = 100102456
$ws.Range("K134").Value2 = 4971.857400000001
$ws.Range("L134").Value2 = 300307368
$ws.Range("M134").Value2 = -2436.857400000001
$ws.Range("N134").Value2 = -300312438

# ============ Sheet: CRP ============
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value2 = 1271.3334
$ws.Range("I16").Value2 = 1512.625
$ws.Range("J16").Value2 = 1078.3
$ws.Range("K16").Value2 = 1512.625
$ws.Range("L16").Value2 = 1078.3
$ws.Range("M16").Value2 = -1225.625
$ws.Range("N16").Value2 = -1652.3
# Row 22
$ws.Range("H22").Value2 = 830
$ws.Range("I22").Value2 = 1236.9
$ws.Range("J22").Value2 = 321.375
$ws.Range("K22").Value2 = 1236.9
$ws.Range("L22").Value2 = 321.375
$ws.Range("M22").Value2 = -886.9000000000001
$ws.Range("N22").Value2 = -1021.375
# Row 113
$ws.Range("H113").Value2 = 1271.3334
$ws.Range("I113").Value2 = 1512.625
$ws.Range("J113").Value2 = 1078.3
$ws.Range("K113").Value2 = 1512.625
$ws.Range("L113").Value2 = 1078.3
$ws.Range("M113").Value2 = 657.375
$ws.Range("N113").Value2 = -5418.3
# Row 132
$ws.Range("H132").Value2 = 31024.205
$ws.Range("I132").Value2 = 38116.15
$ws.Range("K132").Value2 = 114348.45
$ws.Range("M132").Value2 = -111818.45

# ============ Sheet: GSM ============
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value2 = 777.2069
$ws.Range("I97").Value2 = 749.5
$ws.Range("J97").Value2 = 838.7778
$ws.Range("K97").Value2 = 749.5
$ws.Range("L97").Value2 = 838.7778
$ws.Range("M97").Value2 = -253.5
$ws.Range("N97").Value2 = -1830.7778
# Row 113
$ws.Range("H113").Value2 = 1540
$ws.Range("I113").Value2 = 1726.6666
$ws.Range("J113").Value2 = 980
$ws.Range("K113").Value2 = 1726.6666
$ws.Range("L113").Value2 = 980
$ws.Range("M113").Value2 = 443.3334
$ws.Range("N113").Value2 = -5320

# ============ Sheet: LTW ============
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value2 = 993.8261
$ws.Range("I16").Value2 = 862.8
$ws.Range("J16").Value2 = 1239.5
$ws.Range("K16").Value2 = 862.8
$ws.Range("L16").Value2 = 1239.5
$ws.Range("M16").Value2 = -692.8
$ws.Range("N16").Value2 = -1579.5
# Row 133
$ws.Range("H133").Value2 = 44141.25
$ws.Range("J133").Value2 = 44141.25
$ws.Range("L133").Value2 = 44141.25
$ws.Range("N133").Value2 = -49201.25

# ============ Sheet: WVR ============
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value2 = 1529.8125
$ws.Range("I100").Value2 = 600
$ws.Range("J100").Value2 = 1662.6428
$ws.Range("K100").Value2 = 1200
$ws.Range("L100").Value2 = 3325.2856
$ws.Range("M100").Value2 = -659
$ws.Range("N100").Value2 = -4407.2856
# Row 107
$ws.Range("H107").Value2 = 274.54544
$ws.Range("I107").Value2 = 297.18182
$ws.Range("J107").Value2 = 251.90909
$ws.Range("K107").Value2 = 891.54546
$ws.Range("L107").Value2 = 755.72727
$ws.Range("M107").Value2 = 1028.45454
$ws.Range("N107").Value2 = -4595.72727
# Row 132
$ws.Range("H132").Value2 = 2637.2036
$ws.Range("I132").Value2 = 458.82977
$ws.Range("J132").Value2 = 17263.428
$ws.Range("K132").Value2 = 1376.48931
$ws.Range("L132").Value2 = 51790.284
$ws.Range("M132").Value2 = 1153.51069
$ws.Range("N132").Value2 = -56850.284
# Row 136
$ws.Range("H136").Value2 = 2145594.2
$ws.Range("I136").Value2 = 2750408.5
$ws.Range("K136").Value2 = 8251225.5
$ws.Range("M136").Value2 = -8248675.5
